$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet had a duplicated data row (row 3 was an exact copy of row 2).
# Remove the duplicate row entirely, shifting everything up.
$ws.Rows.Item(3).Delete()

# Update the remaining data row with the corrected alert-message sample
# values (leave type, dates re-formatted, reason).
$ws.Cells.Item(2, 1).Value = "Casual Leave"
$ws.Cells.Item(2, 2).Value = "Full day"
$ws.Cells.Item(2, 3).Value = "28-August-2023"
$ws.Cells.Item(2, 4).Value = "30-August-2023"
$ws.Cells.Item(2, 5).Value = "time pass"

# The row previously used an oversized blue "Courier New" font with a
# taller row height - put it back to the plain default style/height.
$ws.Rows.Item(2).Style = "Normal"
$ws.Rows.Item(2).AutoFit()

# The two date columns keep a (plain) text number format.
$ws.Range("C2:D2").NumberFormat = "@"

# Resize the columns to fit their new, shorter content.
$ws.Columns.Item(1).ColumnWidth = 14.5853
$ws.Columns.Item(2).ColumnWidth = 12.7509
$ws.Columns.Item(3).ColumnWidth = 13.4202
$ws.Columns.Item(4).ColumnWidth = 12.5858
$ws.Columns.Item(5).ColumnWidth = 8.5858

# Update the selection/active cell shown when the workbook is reopened.
$ws.Range("D6").Select() | Out-Null

# Restore the (enlarged) window size recorded for the workbook.
$wb.Windows.Item(1).Width = 23370
$wb.Windows.Item(1).Height = 10740 | Out-Null
